# Novas tarefas para entrega
# Insert a brand-new first sheet "NOVO - Tarefas" with the team's task list.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() with no args inserts the new sheet before the currently
# active sheet (the first tab), which is exactly where it needs to land.
$ws = $wb.Worksheets.Add()
$ws.Name = "NOVO - Tarefas"

$ws.Range("A1").Value = "Cadastrar novo perfil (usuário)"
$ws.Range("A2").Value = "Cadastrar novo endereço (usuário)"
$ws.Range("A3").Value = "cadastrar nova massa (administrador)"
$ws.Range("A4").Value = "cadastrar novo recheio (administrador)"
$ws.Range("A5").Value = "cadastrar nova cobertura (administrador)"
$ws.Range("A6").Value = "cadastrar novo cupcake_pedido (adicionar ao carrinho)"
$ws.Range("A7").Value = "cadastrar novo pedido (finalizar compra com varios cupcake_pedido)"
$ws.Range("A8").Value = "Importar arquivo de novo endereço (.csv) + possibilidade de editar dados antes de gravar"
$ws.Range("A9").Value = "Importar arquivo de novo endereço (.xml) + possibilidade de editar dados antes de gravar"
$ws.Range("A10").Value = "Relatório ADM: ver status de todos os pedidos feitos até o momento"
$ws.Range("A11").Value = "Relatório ADM: ver massa mais comprada até hoje"
$ws.Range("A12").Value = "Relatório ADM: ver recheio mais comprado até hoje"
$ws.Range("A13").Value = "Relatório ADM: ver cobertura mais comprada até hoje"
$ws.Range("A14").Value = "Relatório USU: Ver histórico de pedidos"
$ws.Range("A16").Value = "Adicional:"
$ws.Range("A17").Value = "Scripts sql para criar banco e inserir dados"

# Column A sized to comfortably fit the longest line (best-fit column width).
$ws.Columns.Item(1).ColumnWidth = 82.14

# Matches the page margins (metric/A4 defaults) used on the other sheets.
$ws.PageSetup.LeftMargin = 0.511811024 * 72
$ws.PageSetup.RightMargin = 0.511811024 * 72
$ws.PageSetup.TopMargin = 0.78740157499999996 * 72
$ws.PageSetup.BottomMargin = 0.78740157499999996 * 72
$ws.PageSetup.HeaderMargin = 0.31496062000000002 * 72
$ws.PageSetup.FooterMargin = 0.31496062000000002 * 72

# Matches the original author's final cursor position/selection on the sheet.
$ws.Range("A10").Select() | Out-Null
